# Improve api endpoint validation / update dependencies / optimize queries
# -- appends one new telemetry row to each of the four worksheets, mirroring
#    the format of the row immediately above it.
#
# NOTE: this runtime's PS interpreter only reliably binds POSITIONAL
# arguments to function parameters (named "-Param value" args come through
# empty), so every call below passes arguments positionally, in parameter
# order.

$wb = $excel.ActiveWorkbook

function Add-TelemetryRow {
    param(
        $ws,
        [int]$row,
        [double]$timeVal,
        [string]$bVal,
        [string]$cVal,
        [string]$dVal,
        [string]$eVal,
        [int]$fVal,
        $gVal,
        [bool]$gIsText,
        [int]$hVal,
        [int]$iVal
    )

    $prevRow = $row - 1

    # Column A: date/time serial, copy the number format from the cell above
    # so it keeps the custom "YYYY-MM-DD HH:MM:SS" style instead of defaulting
    # to General.
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat
    $aCell.Value = $timeVal

    # Columns B-E: short hex-ish strings, always plain text.
    $ws.Cells.Item($row, 2).Value = $bVal
    $ws.Cells.Item($row, 3).Value = $cVal
    $ws.Cells.Item($row, 4).Value = $dVal
    $ws.Cells.Item($row, 5).Value = $eVal

    # Column F: plain integer.
    $ws.Cells.Item($row, 6).Value = $fVal

    # Column G: usually numeric, but one row needs to stay text because the
    # literal has more significant digits than a double can round-trip.
    $gCell = $ws.Cells.Item($row, 7)
    if ($gIsText) {
        $gCell.NumberFormat = "@"
        $gCell.Value = $gVal
    } else {
        $gCell.Value = $gVal
    }

    # Columns H-I: plain integers.
    $ws.Cells.Item($row, 8).Value = $hVal
    $ws.Cells.Item($row, 9).Value = $iVal
}

# ---- Sheet 1: ROW50-FE-LIFTER -- append row 22 ----
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$g1 = $ws1.Cells.Item(21, 7).Value2
Add-TelemetryRow `
    $ws1 `
    22 `
    45735.12566611111 `
    "0x01,0x90" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," `
    "0x01,0x7e" `
    "0xe" `
    400 `
    $g1 $false `
    382 `
    14

# ---- Sheet 2: ROW50-MID-LIFTER -- append row 24 ----
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
Add-TelemetryRow `
    $ws2 `
    24 `
    45735.10394675926 `
    "0x01,0x90 " `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," `
    "0x01,0x82" `
    "0x19" `
    400 `
    "568631262647113771663628" $true `
    386 `
    25

# ---- Sheet 3: ROW11-FE-LIFTER -- append row 22 ----
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$g3 = $ws3.Cells.Item(21, 7).Value2
Add-TelemetryRow `
    $ws3 `
    22 `
    45735.14782016203 `
    "0x01,0x90" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," `
    "0x01,0x7e" `
    "0x14" `
    400 `
    $g3 $false `
    382 `
    20

# ---- Sheet 4: ROW11-MID-LIFTER -- append row 22 ----
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$g4 = $ws4.Cells.Item(21, 7).Value2
Add-TelemetryRow `
    $ws4 `
    22 `
    45735.29618920139 `
    "0x01,0x90" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," `
    "0x01,0x86" `
    "0x19" `
    400 `
    $g4 $false `
    390 `
    25

Write-Output "done"
